# Add a new credits entry (row 4) mirroring the existing "What / Name / Link"
# rows already on the sheet - a character-base asset credit for Bagong Games.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Character base (dummy)"
$ws.Range("B4").Value = "Bagong Games"
$ws.Range("C4").Value = "https://bagong-games.itch.io/hana-caraka-base-character"

# Resize columns to fit the new (longer) content, like Excel does automatically.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# Leave the selection where the author's session ended up.
$ws.Range("C8").Select()
